# The dataset had a duplicated block: the record that was in row 304
# (Rep 2 / Plot 52 for genotype N19-0538, etc.) was an accidental repeat of
# row 303. Remove the stray row so every subsequent reading shifts up one
# row and realigns with its correct Rep/Plot/lsmeans values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(304).Delete()
